$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-30 Thursday" "2025-01-31 Friday"

Replace-Text "892×7=6244" "456×7=3192"
Replace-Text "319×2=638" "381×9=3429"
Replace-Text "184×7=1288" "635×8=5080"
Replace-Text "189×2=378" "247×9=2223"
Replace-Text "321×6=1926" "291×9=2619"

Replace-Text "978×6=5868" "265×4=1060"
Replace-Text "107×6=642" "644×7=4508"
Replace-Text "866×7=6062" "841×4=3364"
Replace-Text "686×6=4116" "196×3=588"
Replace-Text "264×9=2376" "314×3=942"

Replace-Text "162×3=486" "248×3=744"
Replace-Text "472×5=2360" "126×4=504"
Replace-Text "340×5=1700" "207×5=1035"
Replace-Text "101×9=909" "225×4=900"
Replace-Text "114×8=912" "178×5=890"

Replace-Text "793×7=5551" "429×6=2574"
Replace-Text "161×3=483" "869×2=1738"
Replace-Text "520×2=1040" "629×4=2516"
Replace-Text "684×2=1368" "654×8=5232"
Replace-Text "593×9=5337" "124×7=868"

Replace-Text "296×5=1480" "620×2=1240"
Replace-Text "859×2=1718" "744×9=6696"
Replace-Text "172×3=516" "886×3=2658"
Replace-Text "566×3=1698" "813×2=1626"
Replace-Text "243×4=972" "108×2=216"
